$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Justin (row 16) and Shaiju (row 25) paid their July,18 subscription in full.
$ws.Range("G16").Value = 500
$ws.Range("G25").Value = 500

# Reflect the editor's final selection/scroll position when the file was saved.
$ws.Range("G25").Select()
